$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.839.12"
$ws.Range("E2").Value = "  +7.61%  "

$ws.Range("D3").Value = "1.761.43"
$ws.Range("E3").Value = "  +5.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3820"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.227"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07674"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.466"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.075"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.25%  "

$ws.Range("D16").Value = "1.763.67"
$ws.Range("E16").Value = "  +6.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001151"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.09%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9995"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06783"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "87.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.498"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.82%  "

$ws.Range("D24").Value = "25.817.39"
$ws.Range("E24").Value = "  +7.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.426"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.906"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.49%  "

$ws.Range("D29").Value = "1.959.78"
$ws.Range("E29").Value = "  +6.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.210"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +23.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.148"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.216"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "14.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.800"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08753"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.704"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06750"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.36%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.364"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.50%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02494"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2258"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.291"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6579"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6351"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.899"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.169"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07510"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.90%  "
